$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading paragraph.
# ------------------------------------------------------------------
$metaSearch = $d.Content
$metaSearch.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metaPara = $metaSearch.Paragraphs(1)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new paragraph, right before the final paragraph (the one
#    that holds the italic AI-image prompt text), containing the bold
#    heading text "Play Black and White Slot Game for Free - Review".
#    We insert it *after* the paragraph that precedes the final one
#    (rather than *before* the final, italic, paragraph) so the new
#    paragraph does not inherit any italic character formatting.
# ------------------------------------------------------------------
$imgPromptSearch = $d.Content
$imgPromptSearch.Find.Execute("Create a fun cartoon image of a Maya warrior", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$imgPromptPara = $imgPromptSearch.Paragraphs(1)
$precedingPara = $imgPromptPara.Previous()
$precedingPara.Range.InsertParagraphAfter()

$newPara = $precedingPara.Next()
$newPara.Style = "Normal"
$newParaStart = $newPara.Range.Start

# Type the heading text (inserting at a collapsed point is safe; only
# *formatting* a collapsed/zero-length range is not, in this runtime).
$insertPoint = $d.Range($newParaStart, $newParaStart)
$insertPoint.InsertAfter("Play Black and White Slot Game for Free - Review")

$newParaEnd = $newPara.Range.End - 1
$headingRange = $d.Range($newParaStart, $newParaEnd)
$headingRange.Bold = 1

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph: the old AI
#    image-prompt text becomes the former meta-description sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a fun cartoon image of a Maya warrior wearing glasses, with a big smile on their face. The warrior should be holding a slot machine handle in one hand, and surrounded by colorful flowers and butterflies. The background should feature a jungle landscape with a Mayan pyramid in the distance. Use a mix of bold colors to make the image pop and convey a sense of excitement and fun. The image should be eye-catching and playful, inviting players to explore the game and enjoy their gambling experience.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Want to know about the Black and White slot game? Read our review and play it for free. Learn about its unique features, gameplay, and weaknesses.",
    2
) | Out-Null
